$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '59.116.29'
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).Value = '  +0.84%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.638.26'
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +3.09%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '523.23'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +2.55%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '146.52'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.09%  '
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$ws.Cells.Item(8, 5).Value = '  +0.65%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '2.655.10'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +3.27%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '6.32'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +1.32%  '
$ws.Cells.Item(11, 5).Value = '  +1.73%  '
$ws.Cells.Item(12, 5).Value = '  +0.76%  '
$ws.Cells.Item(13, 5).Value = '  -1.05%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '3.102.47'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +3.26%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '59.099.26'
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +0.93%  '
$ws.Cells.Item(16, 5).Value = '  +0.09%  '
$ws.Cells.Item(17, 5).Value = '  +0.72%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '2.647.63'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  +3.14%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '347.64'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +0.05%  '
$ws.Cells.Item(20, 5).Value = '  -0.68%  '
$ws.Cells.Item(21, 5).Value = '  +0.56%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.16'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +2.72%  '
$ws.Cells.Item(23, 5).Value = '  -0.13%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '61.97'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +2.19%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.418'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +1.16%  '
$ws.Cells.Item(26, 5).Value = '  +3.99%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.994'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -0.29%  '
$ws.Cells.Item(28, 5).Value = '  +0.75%  '
$ws.Cells.Item(29, 5).Value = '  +1.54%  '
$ws.Cells.Item(30, 5).Value = '  -0.07%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.27'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +5.03%  '
$ws.Cells.Item(32, 5).Value = '  +2.97%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '18.99'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +1.91%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '150.21'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +0.40%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.981'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +8.93%  '
$ws.Cells.Item(36, 5).Value = '  +1.33%  '
$ws.Cells.Item(37, 5).Value = '  +1.29%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '36.74'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +2.04%  '
$ws.Cells.Item(39, 5).Value = '  +0.67%  '
$ws.Cells.Item(40, 5).Value = '  +3.47%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.43'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +1.47%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '279.28'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -2.83%  '
$ws.Cells.Item(43, 2).Value = 'Mantle'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.611'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +0.94%  '
$ws.Cells.Item(44, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.995'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -0.06%  '
$ws.Cells.Item(45, 5).Value = '  -1.11%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '19.60'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +2.69%  '
$ws.Cells.Item(47, 5).Value = '  -2.06%  '
$ws.Cells.Item(48, 5).Value = '  +0.43%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.991.83'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +3.49%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '4.68'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +1.10%  '
